$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their original text (General) format so values
# are written back as plain text strings, matching the inline-string cells
# in the source workbook (prevents Excel auto-converting numeric-looking
# text like "354.96" or "1.00" into real numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.932.41'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.831.23'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.44%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '354.96'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +6.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '114.04'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.80'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0854'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.07'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.76'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.277.18'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.821.47'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.899'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.931.95'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.38'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +7.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.15'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.55'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '270.82'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.87'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +5.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.76'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.34'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.31%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0457'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +29.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.79'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '33.93'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.52%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0832'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.91'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.58%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.05%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.96%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.18'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.57%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.093.46'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.72'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.946'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +8.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '60.89'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.40%  '
